$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.172.79"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "'1.849.48"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'234.95"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4691"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").Value = "'0.2887"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("D9").Value = "'0.06554"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "'21.83"
$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("D11").Value = "'0.07950"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "'97.43"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "'1.851.88"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").Value = "'5.087"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").Value = "'0.6742"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").Value = "'268.34"
$ws.Range("E16").Value = "  -3.99%  "

$ws.Range("D17").Value = "'30.133.38"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").Value = "'13.61"
$ws.Range("E18").Value = "  +7.06%  "

$ws.Range("D19").Value = "'0.000007631"
$ws.Range("E19").Value = "  +4.62%  "

$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'2.092.07"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'5.188"
$ws.Range("E23").Value = "  -5.76%  "

$ws.Range("D24").Value = "'6.127"
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("D25").Value = "'166.24"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").Value = "'9.140"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("D27").Value = "'18.78"
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("D28").Value = "'1.926"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("D29").Value = "'1.379"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").Value = "'0.09819"
$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").Value = "'4.276"
$ws.Range("E32").Value = "  -2.37%  "

$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("D34").Value = "'0.04691"
$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("D36").Value = "'0.6971"
$ws.Range("E36").Value = "  -1.35%  "

$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").Value = "'0.01864"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Value = "'2.607"
$ws.Range("E39").Value = "  +3.05%  "

$ws.Range("D40").Value = "'6.317"
$ws.Range("E40").Value = "  +0.57%  "

$ws.Range("D41").Value = "'72.88"
$ws.Range("E41").Value = "  -1.22%  "

$ws.Range("D42").Value = "'1.928"
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").Value = "'0.9987"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").Value = "'0.8364"
$ws.Range("E44").Value = "  -1.38%  "

$ws.Range("D45").Value = "'102.92"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.178"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'938.28"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "'7.000"
$ws.Range("E49").Value = "  -2.52%  "

$ws.Range("D50").Value = "'33.79"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").Value = "'0.05656"
$ws.Range("E51").Value = "  +0.33%  "
